# Updated cryptos list on Sat Nov  9 10:47:22 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.473.50"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "3.042.99"
$ws.Range("E3").Value = "  +4.35%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "202.45"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "626.40"
$ws.Range("E6").Value = "  +4.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.552"
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("E9").Value = "  +5.95%  "
$ws.Range("D10").Value = "3.042.97"
$ws.Range("E10").Value = "  +4.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.439"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.12"
$ws.Range("E13").Value = "  +4.75%  "
$ws.Range("D14").Value = "3.607.65"
$ws.Range("E14").Value = "  +4.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.58"
$ws.Range("E15").Value = "  +6.01%  "
$ws.Range("D16").Value = "76.387.23"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").Value = "3.038.37"
$ws.Range("E18").Value = "  +4.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.52"
$ws.Range("E19").Value = "  +4.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.02"
$ws.Range("E20").Value = "  +1.96%  "
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.31"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.36"
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("E24").Value = "  +2.98%  "
$ws.Range("D25").Value = "3.202.90"
$ws.Range("E25").Value = "  +4.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.40"
$ws.Range("E26").Value = "  +3.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +2.36%  "
$ws.Range("E29").Value = "  +2.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.995"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.29"
$ws.Range("E31").Value = "  +7.44%  "
$ws.Range("E32").Value = "  +1.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "512.38"
$ws.Range("E33").Value = "  +1.94%  "
$ws.Range("E34").Value = "  +7.24%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  +3.52%  "
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.910"
$ws.Range("E37").Value = "  +39.23%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.29"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.387"
$ws.Range("E39").Value = "  +7.38%  "
$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "20.04"
$ws.Range("E40").Value = "  +2.12%  "
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.106"
$ws.Range("E41").Value = "  +1.20%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "189.10"
$ws.Range("E42").Value = "  +4.13%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.113"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.20"
$ws.Range("E45").Value = "  +4.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.29"
$ws.Range("E46").Value = "  +7.51%  "
$ws.Range("E47").Value = "  +5.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.67"
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("E49").Value = "  +4.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.609"
$ws.Range("E50").Value = "  +6.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.89"
$ws.Range("E51").Value = "  +4.72%  "
